# Edit water bill model and query water bill
# Adds a new "Người tạo" column (J) to the header row, matching the
# formatting of the existing header cells, and updates the active
# selection as left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (I1) onto the
# new header cell (J1), then set its text.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("J1").Value = "Người tạo"

# Leave the selection where the author left it before saving.
$ws.Range("J7").Select()
